$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$target = $ws.Range("R4:R5,J4:P5,A8:H9,J6:R9,A4:A7,C4:H7,A10:R95")
Write-Host "Areas count: $($target.Areas.Count())"
for ($i=1; $i -le $target.Areas.Count(); $i++) {
  Write-Host "Area $i : $($target.Areas.Item($i).Address())"
}
